$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Menadžeri imaju mogućnost unosa podataka o zaposlenicima i pregleda njihovih radnih aktivnosti kako bi efikasno koordinirali projekte i timove.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Menadžeri imaju mogućnost unosa podataka o zaposlenicima i pregleda njihovih radnih aktivnosti kako bi efikasno koordinirali projekte i timove, te unos i pregled projekata.",
    2
)
